$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 7486
$ws.Range("I10").Value = 7486
$ws.Range("K10").Value = 7486
$ws.Range("M10").Value = -7193

$ws.Range("H11").Value = 454732.9
$ws.Range("I11").Value = 454732.9
$ws.Range("K11").Value = 454732.9
$ws.Range("M11").Value = -454592.9

$ws.Range("H19").Value = 741.2727
$ws.Range("I19").Value = 403
$ws.Range("J19").Value = 868.125
$ws.Range("K19").Value = 403
$ws.Range("L19").Value = 868.125
$ws.Range("M19").Value = -228
$ws.Range("N19").Value = -1218.125

$ws.Range("H32").Value = 1072
$ws.Range("I32").Value = 890
$ws.Range("K32").Value = 890
$ws.Range("M32").Value = -564

$ws.Range("H132").Value = 2298.681
$ws.Range("I132").Value = 1734.1777
$ws.Range("K132").Value = 5202.5331
$ws.Range("M132").Value = -2672.5331

$ws.Range("H133").Value = 124983.336
$ws.Range("J133").Value = 124983.336
$ws.Range("L133").Value = 124983.336
$ws.Range("N133").Value = -135103.336

$ws.Range("H136").Value = 89979.664
$ws.Range("J136").Value = 89979.664
$ws.Range("L136").Value = 89979.664
$ws.Range("N136").Value = -100179.664

$ws.Range("H137").Value = 1978.2354
$ws.Range("I137").Value = 1873.4166
$ws.Range("K137").Value = 5620.2498
$ws.Range("M137").Value = -3070.2498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 565.5806
$ws.Range("I2").Value = 565.4167
$ws.Range("K2").Value = 565.4167
$ws.Range("M2").Value = -452.4167

$ws.Range("H32").Value = 8613.379999999999
$ws.Range("I32").Value = 4861.213
$ws.Range("K32").Value = 4861.213
$ws.Range("M32").Value = -4574.213

$ws.Range("H63").Value = 3439.8
$ws.Range("J63").Value = 3733.3333
$ws.Range("L63").Value = 3733.3333
$ws.Range("N63").Value = -5105.3333

$ws.Range("H66").Value = 3439.8
$ws.Range("J66").Value = 3733.3333
$ws.Range("L66").Value = 18666.6665
$ws.Range("N66").Value = -25530.6665

$ws.Range("H110").Value = 1695.125
$ws.Range("I110").Value = 1747.2667
$ws.Range("J110").Value = 913
$ws.Range("K110").Value = 1747.2667
$ws.Range("L110").Value = 913
$ws.Range("M110").Value = 297.7333000000001
$ws.Range("N110").Value = -5003

$ws.Range("H116").Value = 565.5806
$ws.Range("I116").Value = 565.4167
$ws.Range("K116").Value = 565.4167
$ws.Range("M116").Value = 1728.5833

$ws.Range("H132").Value = 2763.074
$ws.Range("I132").Value = 2086.9524
$ws.Range("K132").Value = 6260.8572
$ws.Range("M132").Value = -3730.8572

$ws.Range("H134").Value = 95500
$ws.Range("J134").Value = 95500
$ws.Range("L134").Value = 95500
$ws.Range("N134").Value = -105640

$ws.Range("H141").Value = 71479
$ws.Range("I141").Value = 34999.5
$ws.Range("J141").Value = 89718.75
$ws.Range("K141").Value = 34999.5
$ws.Range("L141").Value = 89718.75
$ws.Range("M141").Value = -29819.5
$ws.Range("N141").Value = -100078.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 565.5806
$ws.Range("I3").Value = 565.4167
$ws.Range("K3").Value = 565.4167
$ws.Range("M3").Value = -451.4167

$ws.Range("H105").Value = 4009079.5
$ws.Range("I105").Value = 5008099.5
$ws.Range("K105").Value = 5008099.5
$ws.Range("M105").Value = -5006352.5

$ws.Range("H132").Value = 143299.33
$ws.Range("J132").Value = 143299.33
$ws.Range("L132").Value = 143299.33
$ws.Range("N132").Value = -153419.33

$ws.Range("H134").Value = 3562.0908
$ws.Range("I134").Value = 2874.5715
$ws.Range("K134").Value = 8623.7145
$ws.Range("M134").Value = -6088.7145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8373.4
$ws.Range("I31").Value = 5688.5757
$ws.Range("J31").Value = 13585.117
$ws.Range("K31").Value = 5688.5757
$ws.Range("L31").Value = 13585.117
$ws.Range("M31").Value = -5393.5757
$ws.Range("N31").Value = -14175.117

$ws.Range("H34").Value = 8373.4
$ws.Range("I34").Value = 5688.5757
$ws.Range("J34").Value = 13585.117
$ws.Range("K34").Value = 5688.5757
$ws.Range("L34").Value = 13585.117
$ws.Range("M34").Value = -5486.5757
$ws.Range("N34").Value = -13989.117

$ws.Range("H58").Value = 3055.4119
$ws.Range("I58").Value = 2607.375
$ws.Range("J58").Value = 3453.6667
$ws.Range("K58").Value = 2607.375
$ws.Range("L58").Value = 3453.6667
$ws.Range("M58").Value = -2404.375
$ws.Range("N58").Value = -3859.6667

$ws.Range("H62").Value = 9595.333000000001
$ws.Range("I62").Value = 2797
$ws.Range("J62").Value = 12994.5
$ws.Range("K62").Value = 2797
$ws.Range("L62").Value = 12994.5
$ws.Range("M62").Value = -2173
$ws.Range("N62").Value = -14242.5

$ws.Range("H65").Value = 9595.333000000001
$ws.Range("I65").Value = 2797
$ws.Range("J65").Value = 12994.5
$ws.Range("K65").Value = 13985
$ws.Range("L65").Value = 64972.5
$ws.Range("M65").Value = -10865
$ws.Range("N65").Value = -71212.5

$ws.Range("H86").Value = 2671682.5
$ws.Range("I86").Value = 4449108.5
$ws.Range("J86").Value = 5544
$ws.Range("K86").Value = 4449108.5
$ws.Range("L86").Value = 5544
$ws.Range("M86").Value = -4447985.5
$ws.Range("N86").Value = -7790

$ws.Range("H89").Value = 2671682.5
$ws.Range("I89").Value = 4449108.5
$ws.Range("J89").Value = 5544
$ws.Range("K89").Value = 22245542.5
$ws.Range("L89").Value = 27720
$ws.Range("M89").Value = -22239926.5
$ws.Range("N89").Value = -38952

$ws.Range("H132").Value = 4457.171
$ws.Range("I132").Value = 2691.8484
$ws.Range("K132").Value = 8075.5452
$ws.Range("M132").Value = -5545.5452

$ws.Range("H136").Value = 3055.4119
$ws.Range("I136").Value = 2607.375
$ws.Range("J136").Value = 3453.6667
$ws.Range("K136").Value = 7822.125
$ws.Range("L136").Value = 10361.0001
$ws.Range("M136").Value = -5272.125
$ws.Range("N136").Value = -15461.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 939.0833
$ws.Range("J5").Value = 1287.5
$ws.Range("L5").Value = 3862.5
$ws.Range("N5").Value = -4086.5

$ws.Range("H9").Value = 500000
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").Value = $null

$ws.Range("H135").Value = 939.0833
$ws.Range("J135").Value = 1287.5
$ws.Range("L135").Value = 11587.5
$ws.Range("N135").Value = -16657.5

$ws.Range("H137").Value = 5655.8184
$ws.Range("I137").Value = 2899.6667
$ws.Range("K137").Value = 8699.000100000001
$ws.Range("M137").Value = -3599.000100000001

$ws.Range("H140").Value = 1532.9032
$ws.Range("I140").Value = 974.5454999999999
$ws.Range("K140").Value = 2923.6365
$ws.Range("M140").Value = 2256.3635

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2679.6667
$ws.Range("I80").Value = 2900
$ws.Range("K80").Value = 2900
$ws.Range("M80").Value = -1902

$ws.Range("H83").Value = 2679.6667
$ws.Range("I83").Value = 2900
$ws.Range("K83").Value = 14500
$ws.Range("M83").Value = -9508

$ws.Range("H113").Value = 11009
$ws.Range("I113").Value = 3008.2
$ws.Range("K113").Value = 3008.2
$ws.Range("M113").Value = -838.1999999999998

$ws.Range("H122").Value = 5480.8335
$ws.Range("I122").Value = 2940.75
$ws.Range("K122").Value = 8822.25
$ws.Range("M122").Value = -6372.25

$ws.Range("H126").Value = 11339
$ws.Range("I126").Value = 12344.5
$ws.Range("J126").Value = 9998.333000000001
$ws.Range("K126").Value = 37033.5
$ws.Range("L126").Value = 29994.999
$ws.Range("M126").Value = -34563.5
$ws.Range("N126").Value = -34934.999

$ws.Range("H132").Value = 8620.357
$ws.Range("I132").Value = 5915.636
$ws.Range("K132").Value = 17746.908
$ws.Range("M132").Value = -15216.908

$ws.Range("H136").Value = 14762.8125
$ws.Range("J136").Value = 14762.8125
$ws.Range("L136").Value = 44288.4375
$ws.Range("N136").Value = -49388.4375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7412.857
$ws.Range("I132").Value = 5648.3335
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 16945.0005
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -14415.0005
$ws.Range("N132").Value = -59060

$ws.Range("H136").Value = 7406.772
$ws.Range("I136").Value = 5110.276
$ws.Range("K136").Value = 15330.828
$ws.Range("M136").Value = -12780.828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3301.225
$ws.Range("I107").Value = 3409.5833
$ws.Range("J107").Value = 3138.6875
$ws.Range("K107").Value = 10228.7499
$ws.Range("L107").Value = 9416.0625
$ws.Range("M107").Value = -8308.749899999999
$ws.Range("N107").Value = -13256.0625

$ws.Range("H113").Value = 343.82144
$ws.Range("I113").Value = 339.5
$ws.Range("J113").Value = 369.75
$ws.Range("K113").Value = 1018.5
$ws.Range("L113").Value = 1109.25
$ws.Range("M113").Value = 1151.5
$ws.Range("N113").Value = -5449.25

$ws.Range("H132").Value = 4357.24
$ws.Range("I132").Value = 3788.7917
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 11366.3751
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -8836.375100000001
$ws.Range("N132").Value = -59060

$ws.Range("H136").Value = 4399.8613
$ws.Range("I136").Value = 3693.5483
$ws.Range("K136").Value = 11080.6449
$ws.Range("M136").Value = -8530.644899999999
